$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.493.79'
$ws.Range("E2").Value = '  -1.37%  '
$ws.Range("D3").Value = '1.910.65'
$ws.Range("E3").Value = '  -2.09%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''239.43'
$ws.Range("E5").Value = '  -1.22%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = '''0.4774'
$ws.Range("E7").Value = '  -2.06%  '
$ws.Range("D8").Value = '''0.2841'
$ws.Range("E8").Value = '  -3.39%  '
$ws.Range("D9").Value = '''0.06681'
$ws.Range("E9").Value = '  -3.59%  '
$ws.Range("E10").Value = '  -4.38%  '
$ws.Range("D11").Value = '''100.92'
$ws.Range("E11").Value = '  -5.68%  '
$ws.Range("D12").Value = '1.912.93'
$ws.Range("E12").Value = '  -2.95%  '
$ws.Range("D13").Value = '''0.07671'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("D15").Value = '''0.6675'
$ws.Range("E15").Value = '  -3.99%  '
$ws.Range("D16").Value = '30.499.33'
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("D17").Value = '''253.25'
$ws.Range("E17").Value = '  -9.06%  '
$ws.Range("D18").Value = '''1.000'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '''0.000007454'
$ws.Range("E19").Value = '  -3.68%  '
$ws.Range("E20").Value = '  -4.08%  '
$ws.Range("D21").Value = '''5.376'
$ws.Range("E21").Value = '  -1.75%  '
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("E23").Value = '  -3.09%  '
$ws.Range("D24").Value = '''9.315'
$ws.Range("E24").Value = '  -4.31%  '
$ws.Range("D25").Value = '''167.06'
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("D26").Value = '''19.01'
$ws.Range("E26").Value = '  -3.43%  '
$ws.Range("D27").Value = '''2.054'
$ws.Range("E27").Value = '  -5.21%  '
$ws.Range("D28").Value = '''4.738'
$ws.Range("E28").Value = '  +3.59%  '
$ws.Range("D29").Value = '''0.1008'
$ws.Range("E29").Value = '  -3.08%  '
$ws.Range("E30").Value = '  -1.41%  '
$ws.Range("D31").Value = '''1.512'
$ws.Range("E31").Value = '  -2.78%  '
$ws.Range("D32").Value = '''4.245'
$ws.Range("E32").Value = '  -3.03%  '
$ws.Range("D33").Value = '''0.04718'
$ws.Range("E33").Value = '  -3.07%  '
$ws.Range("D34").Value = '''0.7271'
$ws.Range("E34").Value = '  -3.26%  '
$ws.Range("D35").Value = '''1.107'
$ws.Range("E35").Value = '  -4.73%  '
$ws.Range("D36").Value = '''0.9999'
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").Value = '''2.706'
$ws.Range("E37").Value = '  -0.97%  '
$ws.Range("D38").Value = '''0.01909'
$ws.Range("E38").Value = '  -4.20%  '
$ws.Range("E39").Value = '  -2.74%  '
$ws.Range("D40").Value = '''74.68'
$ws.Range("E40").Value = '  -3.32%  '
$ws.Range("D41").Value = '''6.203'
$ws.Range("E41").Value = '  -4.47%  '
$ws.Range("D42").Value = '''1.962'
$ws.Range("E42").Value = '  -6.42%  '
$ws.Range("D43").Value = '''0.8612'
$ws.Range("E43").Value = '  -4.12%  '
$ws.Range("D44").Value = '''104.98'
$ws.Range("D45").Value = '''0.9999'
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").Value = '''0.4229'
$ws.Range("E46").Value = '  -4.40%  '
$ws.Range("D47").Value = '''7.396'
$ws.Range("D48").Value = '''930.39'
$ws.Range("E48").Value = '  -6.48%  '
$ws.Range("D49").Value = '''0.1198'
$ws.Range("E49").Value = '  -3.82%  '
$ws.Range("D50").Value = '''34.72'
$ws.Range("E50").Value = '  -2.85%  '
$ws.Range("D51").Value = '''8.781'
$ws.Range("E51").Value = '  -4.20%  '
